# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the two obsolete employee records (GABRIEL JOSE NOVOA PAUTT
#    in row 16 and HARLIN SILVA CASTRO in row 17). Deleting the whole
#    rows shifts everything below up by two and keeps the per-row
#    borders/styles intact.
# ------------------------------------------------------------------
$ws.Range("A17:J17").EntireRow.Delete() | Out-Null
$ws.Range("A16:J16").EntireRow.Delete() | Out-Null

# ------------------------------------------------------------------
# 2. Rewrite the remaining 6 worker rows (16-21) with the new data set
#    and ordering used in the updated statement of account.
# ------------------------------------------------------------------
$data = @(
    @("CC", "1096207112", "INYERMAN JOSE FRANCO ALVARADO", "2410", 2133, 1600000),
    @("CC", "8703326", "JOSE ANTONIO TORREGROSA MIRANDA", "2410", 23467, 1600000),
    @("CC", "1143346189", "NAYADETH DEL CARMEN POMARES CASTELLAR", "2505", 56940, 1423500),
    @("CC", "1143389181", "MARIA ELENA PEREZ ORTEGA", "2505", 56940, 1423500),
    @("CC", "1047405911", "SAIRO LUIS DE AVILA MEDRANO", "2505", 56940, 1423500),
    @("CC", "1043298004", "SARAY SOFIA ARRIETA JIMENEZ", "2505", 64000, 1600000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row++
}

# ------------------------------------------------------------------
# 3. Update the summary figures at the top of the statement.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 260420      # VALOR MORA total
$ws.Range("C13").Value = 6           # Cant. Trabajadores
$ws.Range("F13").Value = 2           # Cant. Periodos
